# Sample Project / Main.xlsx - "Rules" sheet
# Cell B11 currently shows the rule id "R40"; the saved project now stores
# "1" there instead (same position in the decision table, same formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

# Stash the cell's current formatting (border/fill/etc.) in an unused cell so
# we can re-apply it after the value swap - entering "1" naively would be
# auto-detected as a number and pick up a different (numeric) style.
$scratch = $ws.Range("Z100")
$target.Copy()
$scratch.PasteSpecial(-4122)

# Force the new value to be stored as text, same as the original "R40".
$target.NumberFormat = "@"
$target.Value = "1"

# Restore the original cell formatting/style.
$scratch.Copy()
$target.PasteSpecial(-4122)
$scratch.Clear()
